$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

Set-TextValue 'D2' '66.952.64'
Set-TextValue 'E2' '  -2.37%  '
Set-TextValue 'D3' '2.660.07'
Set-TextValue 'E3' '  -1.47%  '
Set-TextValue 'E4' '  -0.01%  '
Set-TextValue 'D5' '592.73'
Set-TextValue 'E5' '  -1.05%  '
Set-TextValue 'D6' '163.48'
Set-TextValue 'E6' '  +2.30%  '
Set-TextValue 'E7' '  +0.00%  '
Set-TextValue 'D8' '0.542'
Set-TextValue 'E8' '  -0.51%  '
Set-TextValue 'D9' '2.659.96'
Set-TextValue 'E9' '  -1.43%  '
Set-TextValue 'D10' '0.139'
Set-TextValue 'E10' '  +0.09%  '
Set-TextValue 'E11' '  +0.84%  '
Set-TextValue 'D12' '0.354'
Set-TextValue 'E12' '  -1.91%  '
Set-TextValue 'D13' '5.17'
Set-TextValue 'E13' '  -2.45%  '
Set-TextValue 'B14' 'WrappedliquidstakedEther2.0'
Set-TextValue 'C14' 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextValue 'D14' '3.148.37'
Set-TextValue 'E14' '  -1.27%  '
Set-TextValue 'B15' 'Avalanche'
Set-TextValue 'C15' 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextValue 'D15' '27.46'
Set-TextValue 'E15' '  -3.05%  '
Set-TextValue 'D16' '0.0000181'
Set-TextValue 'E16' '  -3.59%  '
Set-TextValue 'D17' '66.963.35'
Set-TextValue 'E17' '  -2.39%  '
Set-TextValue 'D18' '2.652.09'
Set-TextValue 'E18' '  -1.67%  '
Set-TextValue 'D19' '11.58'
Set-TextValue 'E19' '  -3.17%  '
Set-TextValue 'D20' '359.29'
Set-TextValue 'E20' '  -2.14%  '
Set-TextValue 'D21' '7.46'
Set-TextValue 'E21' '  -3.27%  '
Set-TextValue 'D22' '4.34'
Set-TextValue 'E22' '  -4.62%  '
Set-TextValue 'D23' '4.76'
Set-TextValue 'E23' '  -2.87%  '
Set-TextValue 'D24' '2.01'
Set-TextValue 'E24' '  -5.74%  '
Set-TextValue 'B25' 'Dai'
Set-TextValue 'C25' 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue 'D25' '1.00'
Set-TextValue 'E25' '  +0.10%  '
Set-TextValue 'B26' 'Litecoin'
Set-TextValue 'C26' 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextValue 'D26' '70.91'
Set-TextValue 'E26' '  -4.96%  '
Set-TextValue 'D27' '9.98'
Set-TextValue 'E27' '  -1.12%  '
Set-TextValue 'D28' '2.818.25'
Set-TextValue 'E28' '  -0.63%  '
Set-TextValue 'E29' '  +0.15%  '
Set-TextValue 'D30' '0.0000101'
Set-TextValue 'E30' '  -3.66%  '
Set-TextValue 'D31' '548.29'
Set-TextValue 'E31' '  -4.63%  '
Set-TextValue 'D32' '7.91'
Set-TextValue 'E32' '  -4.53%  '
Set-TextValue 'D33' '1.37'
Set-TextValue 'E33' '  -5.49%  '
Set-TextValue 'D34' '1.91'
Set-TextValue 'E34' '  -1.72%  '
Set-TextValue 'B35' 'FirstDigitalUSD'
Set-TextValue 'C35' 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue 'D35' '0.999'
Set-TextValue 'E35' '  +0.05%  '
Set-TextValue 'B36' 'Kaspa'
Set-TextValue 'C36' 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue 'D36' '0.128'
Set-TextValue 'E36' '  -3.21%  '
Set-TextValue 'D37' '1.55'
Set-TextValue 'E37' '  -6.03%  '
Set-TextValue 'D38' '19.36'
Set-TextValue 'E38' '  -3.36%  '
Set-TextValue 'D39' '154.16'
Set-TextValue 'E39' '  -3.95%  '
Set-TextValue 'D40' '0.370'
Set-TextValue 'E40' '  -2.78%  '
Set-TextValue 'D41' '5.23'
Set-TextValue 'E41' '  -3.89%  '
Set-TextValue 'D42' '1.81'
Set-TextValue 'E42' '  -5.55%  '
Set-TextValue 'D43' '17.89'
Set-TextValue 'E43' '  +0.14%  '
Set-TextValue 'B44' 'USDe'
Set-TextValue 'C44' 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
Set-TextValue 'D44' '1.00'
Set-TextValue 'E44' '  +0.03%  '
Set-TextValue 'B45' 'dogwifhat'
Set-TextValue 'C45' 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue 'D45' '2.49'
Set-TextValue 'E45' '  -6.24%  '
Set-TextValue 'D46' '40.11'
Set-TextValue 'E46' '  -0.92%  '
Set-TextValue 'D47' '0.0₆0294'
Set-TextValue 'E47' '  -7.44%  '
Set-TextValue 'D48' '0.581'
Set-TextValue 'E48' '  -3.11%  '
Set-TextValue 'D49' '151.76'
Set-TextValue 'E49' '  -4.22%  '
Set-TextValue 'D50' '3.79'
Set-TextValue 'E50' '  -3.94%  '
Set-TextValue 'D51' '1.71'
Set-TextValue 'E51' '  -4.03%  '
